$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("angles")

# Update vector component values
$ws.Range("C10").Value = 1
$ws.Range("E10").Value = 1
$ws.Range("E11").Value = 1

# Recalculate so the dependent formula (G10) picks up new cached value
$excel.Calculate()

# Update the active selection on the sheet
$ws.Activate()
$ws.Range("E12").Select()
